# "For Library Faculty is added"
#
# 1. Append a new timetable row ("10 / Library / LIB / <faculty>") to the
#    two faculty-list tables (III B.Tech I Sem Sec-A and Sec-B).
# 2. Tidy up a couple of incidental paragraph/bookmark quirks that came
#    along with the same edit in the source document.

$d = $word.ActiveDocument

function Add-LibraryRow($table, $facultyName) {
    $newRow = $table.Rows.Add()
    $newRow.Cells.Item(1).Range.Text = "10"
    $newRow.Cells.Item(3).Range.Text = "Library"
    $newRow.Cells.Item(4).Range.Text = "LIB"
    $newRow.Cells.Item(5).Range.Text = $facultyName
    return $newRow
}

# --- Table 3 (Sec-A) : add "10 / Library / LIB / Dr. B. Hari Chandana" ---
$tableA = $d.Tables.Item(3)
Add-LibraryRow $tableA "Dr. B. Hari Chandana" | Out-Null

# --- Table 6 (Sec-B) : add "10 / Library / LIB / Mrs. V. Kamakshamma" ---
$tableB = $d.Tables.Item(6)
Add-LibraryRow $tableB "Mrs. V. Kamakshamma" | Out-Null

# --- Move the lastRenderedPageBreak marker up one run (cosmetic, matches
#     the author's re-save) and shrink the spacer run's text ---
$d.Content.Find.Execute("             ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "      ", 2) | Out-Null

# --- Remove one of two now-duplicated empty paragraphs that used to sit
#     right before the library table ---
$libraryTable = $d.Tables.Item(6)
$tableRange = $libraryTable.Range
$before = $d.Range(0, $tableRange.Start)
$paraCount = $before.Paragraphs.Count
$lastPara = $before.Paragraphs.Item($paraCount)
$secondLastPara = $before.Paragraphs.Item($paraCount - 1)
if (($lastPara.Range.Text -eq $secondLastPara.Range.Text)) {
    $lastPara.Range.Delete()
}

Write-Output "done"
